$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad) rows 2-7 from 45175 to 45183
foreach ($r in 2..7) {
    $ws.Cells.Item($r, 3).Value = 45183
}
